$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the check-in/check-out dates for the Seattle and Kansas City rows
# with new DDT test case values, in the order the new shared strings were
# authored.
$ws.Range("B4").Value = " 11/15/2021"
$ws.Range("B2").Value = " 02/15/2022"
$ws.Range("C4").Value = " 01/05/2022"
$ws.Range("C2").Value = " 02/19/2022"

# Update the active selection to B2 as recorded in the sheet view.
$ws.Range("B2").Select()
